$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 6 de Octubre de 2020 a las 00:26"

$ws.Range("B4").Value = 7672810
$ws.Range("C4").Value = 35898
$ws.Range("D4").Value = 4886174
$ws.Range("E4").Value = 2571709
$ws.Range("G4").Value = 316
$ws.Range("H4").Value = 214927
$ws.Range("D6").Value = 4295302
$ws.Range("E6").Value = 485258
$ws.Range("B10").Value = 829999
$ws.Range("C10").Value = 1830
$ws.Range("D10").Value = 712888
$ws.Range("E10").Value = 84277
$ws.Range("G10").Value = 92
$ws.Range("H10").Value = 32834
$ws.Range("B26").Value = 304657
$ws.Range("C26").Value = 3086
$ws.Range("E26").Value = 31341
$ws.Range("B29").Value = 168024
$ws.Range("C29").Value = 1868
$ws.Range("D29").Value = 141687
$ws.Range("E29").Value = 16845
$ws.Range("B41").Value = 103781
$ws.Range("C41").Value = 98
$ws.Range("D41").Value = 97398
$ws.Range("E41").Value = 393
$ws.Range("G41").Value = 9
$ws.Range("H41").Value = 5990
$ws.Range("B57").Value = 73116
$ws.Range("C57").Value = 454
$ws.Range("D57").Value = 67933
$ws.Range("E57").Value = 4922
$ws.Range("B58").Value = 59465
$ws.Range("C58").Value = 120
$ws.Range("D58").Value = 50951
$ws.Range("E58").Value = 7401
$ws.Range("B84").Value = 21870
$ws.Range("C84").Value = 283
$ws.Range("D84").Value = 15179
$ws.Range("E84").Value = 5837
$ws.Range("G84").Value = 10
$ws.Range("H84").Value = 854
$ws.Range("B101").Value = 11654
$ws.Range("C101").Value = 28
$ws.Range("D101").Value = 9451
$ws.Range("E101").Value = 2078
$ws.Range("G101").Value = 2
$ws.Range("H101").Value = 125
$ws.Range("A102").Value = "Guinea"
$ws.Range("B102").Value = 10800
$ws.Range("C102").Value = 46
$ws.Range("D102").Value = 10161
$ws.Range("E102").Value = 572
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 67
$ws.Range("A103").Value = "Consejo Danes para los Refugiados"
$ws.Range("B103").Value = 10778
$ws.Range("C103").Value = 18
$ws.Range("D103").Value = 10239
$ws.Range("E103").Value = 265
$ws.Range("H103").Value = 274
$ws.Range("B112").Value = 8808
$ws.Range("C112").Value = 11
$ws.Range("D112").Value = 8135
$ws.Range("E112").Value = 619
$ws.Range("B114").Value = 7898
$ws.Range("C114").Value = 10
$ws.Range("D114").Value = 6424
$ws.Range("E114").Value = 1246
$ws.Range("B118").Value = 6433
$ws.Range("C118").Value = 73
$ws.Range("D118").Value = 5524
$ws.Range("E118").Value = 841
$ws.Range("G118").Value = 3
$ws.Range("H118").Value = 68
$ws.Range("B121").Value = 5579
$ws.Range("C121").Value = 10
$ws.Range("D121").Value = 5141
$ws.Range("E121").Value = 326
$ws.Range("G121").Value = 1
$ws.Range("H121").Value = 112
$ws.Range("A122").Value = "Angola"
$ws.Range("B122").Value = 5530
$ws.Range("C122").Value = 128
$ws.Range("D122").Value = 2591
$ws.Range("E122").Value = 2740
$ws.Range("G122").Value = 4
$ws.Range("H122").Value = 199
$ws.Range("A123").Value = "Guadalupe"
$ws.Range("B123").Value = 5528
$ws.Range("C123").Value = 0
$ws.Range("D123").Value = 2199
$ws.Range("E123").Value = 3272
$ws.Range("H123").Value = 57
$ws.Range("A124").Value = "Republica de Yibuti"
$ws.Range("B124").Value = 5421
$ws.Range("C124").Value = 2
$ws.Range("D124").Value = 5352
$ws.Range("E124").Value = 8
$ws.Range("H124").Value = 61
$ws.Range("B135").Value = 4411
$ws.Range("C135").Value = 45
$ws.Range("D135").Value = 1168
$ws.Range("E135").Value = 3036
$ws.Range("G135").Value = 2
$ws.Range("H135").Value = 207
$ws.Range("A146").Value = "Guyana"
$ws.Range("B146").Value = 3188
$ws.Range("C146").Value = 95
$ws.Range("D146").Value = 1972
$ws.Range("E146").Value = 1126
$ws.Range("G146").Value = 3
$ws.Range("H146").Value = 90
$ws.Range("A147").Value = "Botsuana"
$ws.Range("B147").Value = 3172
$ws.Range("D147").Value = 710
$ws.Range("E147").Value = 2446
$ws.Range("H147").Value = 16
$ws.Range("B155").Value = 2184
$ws.Range("C155").Value = 17
$ws.Range("D155").Value = 1420
$ws.Range("E155").Value = 705
$ws.Range("B165").Value = 1354
$ws.Range("C165").Value = 6
$ws.Range("E165").Value = 36
$ws.Range("D169").Value = 888
$ws.Range("E169").Value = 10
$ws.Range("D190").Value = 210
$ws.Range("E190").Value = 2
$ws.Range("B191").Value = 200
$ws.Range("C191").Value = 1
$ws.Range("E191").Value = 11
